# Update the Abstract (D7) and Authors (E7) cells for row 7 to the
# cleaned-up text (stray id="ParN"> markers removed, double blank lines
# collapsed, and the author list's separator spacing widened), matching
# the "working on per-file dump / bug fix" commit for reference #72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = 'Objective
To investigate the occurrence of olfactory and gustatory dysfunctions in patients with laboratory-confirmed COVID-19 infection.
Methods
Patients with laboratory-confirmed COVID-19 infection were recruited from 12 European hospitals.
 The following epidemiological and clinical outcomes have been studied: age, sex, ethnicity, comorbidities, and general and otolaryngological symptoms.
 Patients completed olfactory and gustatory questionnaires based on the smell and taste component of the National Health and Nutrition Examination Survey, and the short version of the Questionnaire of Olfactory Disorders-Negative Statements (sQOD-NS).
Results
A total of 417 mild-to-moderate COVID-19 patients completed the study (263 females).
 The most prevalent general symptoms consisted of cough, myalgia, and loss of appetite.
 Face pain and nasal obstruction were the most disease-related otolaryngological symptoms.
 85.6% and 88.0% of patients reported olfactory and gustatory dysfunctions, respectively.
 There was a significant association between both disorders (p &lt; 0.001).
 Olfactory dysfunction (OD) appeared before the other symptoms in 11.8% of cases.
 The sQO-NS scores were significantly lower in patients with anosmia compared with normosmic or hyposmic individuals (p = 0.001).
 Among the 18.2% of patients without nasal obstruction or rhinorrhea, 79.7% were hyposmic or anosmic.
 The early olfactory recovery rate was 44.0%.
 Females were significantly more affected by olfactory and gustatory dysfunctions than males (p = 0.001).
Conclusion
Olfactory and gustatory disorders are prevalent symptoms in European COVID-19 patients, who may not have nasal symptoms.
 The sudden anosmia or ageusia need to be recognized by the international scientific community as important symptoms of the COVID-19 infection.
Electronic supplementary material
The online version of this article (10.1007/s00405-020-05965-1) contains supplementary material, which is available to authorized users.
'
$ws.Range("E7").Value = '[Jerome R.%Lechien%Jerome.Lechien@umons.ac.be%0,      Carlos M.%Chiesa-Estomba%NULL%0,      Daniele R.%De Siati%NULL%0,      Mihaela%Horoi%NULL%0,      Serge D.%Le Bon%NULL%0,      Alexandra%Rodriguez%NULL%0,      Didier%Dequanter%NULL%0,      Serge%Blecic%NULL%0,      Fahd%El Afia%NULL%0,      Lea%Distinguin%NULL%0,      Younes%Chekkoury-Idrissi%NULL%0,      Stéphane%Hans%NULL%0,      Irene Lopez%Delgado%NULL%0,      Christian%Calvo-Henriquez%NULL%0,      Philippe%Lavigne%NULL%0,      Chiara%Falanga%NULL%0,      Maria Rosaria%Barillari%NULL%0,      Giovanni%Cammaroto%NULL%0,      Mohamad%Khalife%NULL%0,      Pierre%Leich%NULL%0,      Christel%Souchay%NULL%0,      Camelia%Rossi%NULL%0,      Fabrice%Journe%NULL%0,      Julien%Hsieh%NULL%0,      Myriam%Edjlali%NULL%0,      Robert%Carlier%NULL%0,      Laurence%Ris%NULL%0,      Andrea%Lovato%NULL%0,      Cosimo%De Filippis%NULL%0,      Frederique%Coppee%NULL%0,      Nicolas%Fakhry%NULL%0,      Tareck%Ayad%NULL%0,      Sven%Saussez%NULL%0]'

# Re-run autofit on the row so that no explicit/custom row height sticks
# around just because the new cell text contains embedded line breaks
# (the source workbook never had an explicit row height here).
$ws.Rows.Item(7).EntireRow.AutoFit()
